$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Price(D), Volume(E)) updates. A leading "~" marker on the price
# value means it is a bare-looking number that must be force-typed as text
# (apostrophe prefix) so Excel keeps it verbatim instead of parsing it as a
# numeric value (which would silently drop trailing/grouping zeros).
$updates = @{
    2 = @("22.150.45", "  -1.45%  ")
    3 = @("1.562.99", "  -0.78%  ")
    4 = @("~1.0000", "  -0.06%  ")
    5 = @($null, "  +0.01%  ")
    6 = @("~289.63", "  +0.54%  ")
    7 = @("~0.3803", "  +3.15%  ")
    8 = @("~0.3293", "  -1.09%  ")
    9 = @("~43.39", "  -9.29%  ")
    10 = @("~1.140", "  -0.94%  ")
    11 = @("~0.07381", "  -2.33%  ")
    12 = @($null, "  -0.06%  ")
    13 = @("~20.02", "  -3.73%  ")
    14 = @("~5.831", "  -1.97%  ")
    15 = @($null, "  -0.71%  ")
    16 = @("1.567.95", "  +0.01%  ")
    17 = @("~0.00001096", "  -2.25%  ")
    18 = @("~0.06662", "  -1.13%  ")
    19 = @("~85.85", "  -2.45%  ")
    20 = @("~6.458", "  +0.99%  ")
    21 = @("~0.9998", "  -0.02%  ")
    22 = @("~16.15", "  -2.39%  ")
    23 = @("~11.71", "  -2.61%  ")
    24 = @("22.163.85", "  -1.31%  ")
    25 = @("~2.270", "  -4.95%  ")
    26 = @("~2.561", "  -3.13%  ")
    27 = @("~150.92", "  +0.28%  ")
    28 = @("~19.13", "  -2.85%  ")
    29 = @("~4.866", "  -2.47%  ")
    30 = @("1.745.87", "  -0.08%  ")
    31 = @("~121.40", "  -3.31%  ")
    32 = @("~1.128", "  +3.36%  ")
    33 = @("~6.051", "  -0.98%  ")
    34 = @("~1.846", "  -7.50%  ")
    35 = @("~9.405", "  -4.65%  ")
    36 = @("~0.08171", "  -2.24%  ")
    37 = @("~5.303", "  -1.34%  ")
    38 = @("~0.06235", "  -2.53%  ")
    39 = @("~0.02302", "  -6.70%  ")
    40 = @($null, "  -4.19%  ")
    41 = @("~1.240", "  -4.04%  ")
    42 = @("~11.09", "  -3.39%  ")
    44 = @("~0.6007", "  -4.27%  ")
    45 = @("~13.82", "  -1.68%  ")
    47 = @("~0.5805", "  -5.13%  ")
    48 = @("~1.982", "  -3.66%  ")
    49 = @("~120.89", "  -3.61%  ")
    50 = @("~1.172", "  -3.36%  ")
    51 = @("~0.06987", "  -3.27%  ")
}

foreach ($row in $updates.Keys) {
    $price = $updates[$row][0]
    $volume = $updates[$row][1]
    if ($price -ne $null) {
        if ($price.StartsWith("~")) {
            # Force text entry so the numeric-looking string (e.g. '1.0000')
            # is preserved exactly, matching how the source value was typed.
            $ws.Cells.Item($row, 4).Value = "'" + $price.Substring(1)
        } else {
            $ws.Cells.Item($row, 4).Value = $price
        }
    }
    $ws.Cells.Item($row, 5).Value = $volume
}
